$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 10276
$ws1.Range("F12").Value = 12841
$ws1.Range("F30").Value = 4284
$ws1.Range("F32").Value = 3839
$ws1.Range("F33").Value = 894
$ws1.Range("F34").Value = 2659
$ws1.Range("F37").Value = 1374
$ws1.Range("F40").Value = 51

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F19").Value = 40

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 10276
$ws4.Range("F11").Value = 12841
$ws4.Range("F29").Value = 4284
$ws4.Range("F30").Value = 3839
$ws4.Range("F31").Value = 894
$ws4.Range("F32").Value = 2659
$ws4.Range("F36").Value = 1374
$ws4.Range("F39").Value = 51
$ws4.Range("F42").Value = 40
